$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35, pushing existing rows 35-100 down to 36-101.
$ws.Rows("35:35").Insert()

# Populate the newly inserted row 35 with the new record's data.
$ws.Range("A35").Value = 7
$ws.Range("B35").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C35").Value = "Ñuble"
$ws.Range("D35").Value = 45251
$ws.Range("E35").Value = 16
$ws.Range("F35").Value = 100112026
$ws.Range("G35").Value = "Haba"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 120
$ws.Range("K35").Value = 9500
$ws.Range("L35").Value = 10000
$ws.Range("M35").Value = 9750
$ws.Range("N35").Value = "$/saco 25 kilos"
$ws.Range("O35").Value = "Provincia de Diguillín"
$ws.Range("P35").Value = 390
$ws.Range("Q35").Value = 25
$ws.Range("R35").Value = "Hortaliza"
